$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("E2").Value = "2024.03.16 13:00-03.17 19:00"
$ws.Range("F2").Value = 889
$ws.Range("E3").Value = "2024.03.16 10:00-03.16 21:00"
$ws.Range("E4").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("F4").Value = 1136
$ws.Range("E5").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F5").Value = 535
$ws.Range("E6").Value = "2024.03.23 10:00-03.24 17:00"
$ws.Range("F6").Value = 230
$ws.Range("E7").Value = "2024.03.23 13:00-03.23 19:00"
$ws.Range("E8").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F8").Value = 692
$ws.Range("E9").Value = "2024.03.23 09:30-03.24 18:00"
$ws.Range("F9").Value = 265
$ws.Range("E10").Value = "2024.03.23 09:30-03.23 18:00"
$ws.Range("E11").Value = "2024.03.24 09:30-03.24 18:00"
$ws.Range("F11").Value = 101
$ws.Range("E12").Value = "2024.03.24 10:00-03.24 17:30"
$ws.Range("F12").Value = 228
$ws.Range("E13").Value = "2024.03.24 11:00-03.24 17:00"
$ws.Range("F13").Value = 165
$ws.Range("E14").Value = "2024.03.30 09:00-03.31 17:00"
$ws.Range("F14").Value = 3349
$ws.Range("E15").Value = "2024.03.30 14:00-03.30 18:00"
$ws.Range("F15").Value = 14
$ws.Range("E16").Value = "2024.03.30 10:00-03.30 17:00"
$ws.Range("F16").Value = 14
$ws.Range("E17").Value = "2024.04.04 10:00-04.05 17:30"
$ws.Range("E18").Value = "2024.04.05 11:30-04.06 18:00"
$ws.Range("F18").Value = 42
$ws.Range("E19").Value = "2024.04.06 10:00-04.06 17:00"
$ws.Range("E20").Value = "2024.04.06 10:00-04.06 17:00"
$ws.Range("F20").Value = 290
$ws.Range("E21").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E22").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E23").Value = "2024.04.14 10:00-04.14 17:00"
$ws.Range("F23").Value = 8
$ws.Range("E24").Value = "2024.04.20 10:00-04.20 17:00"
$ws.Range("F24").Value = 681
$ws.Range("E25").Value = "2024.05.04 10:00-05.04 17:00"
$ws.Range("F25").Value = 62
$ws.Range("E26").Value = "2024.05.04 09:30-05.04 16:00"
$ws.Range("F26").Value = 258
$ws.Range("E27").Value = "2024.05.05 10:00-05.05 17:00"
$ws.Range("F27").Value = 972
$ws.Range("E28").Value = "2024.05.10 10:30-05.12 18:30"
$ws.Range("F28").Value = 56
$ws.Range("E29").Value = "2024.05.18 10:00-05.18 17:00"
$ws.Range("F29").Value = 1614
$ws.Range("E30").Value = "2024.05.18 10:00-05.18 17:00"
$ws.Range("F30").Value = 351
$ws.Range("E31").Value = "2024.05.25 10:30-05.25 17:00"

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("E2").Value = "2024.03.17 19:00-03.17 20:30"
$ws.Range("E3").Value = "2024.03.17 19:30-03.17 21:00"
$ws.Range("F3").Value = 34
$ws.Range("E4").Value = "2024.03.23 20:00-03.23 22:00"
$ws.Range("E5").Value = "2024.03.31 19:00-03.31 20:30"
$ws.Range("F5").Value = 243
$ws.Range("E6").Value = "2024.04.13 19:00-04.13 22:00"
$ws.Range("E7").Value = "2024.04.14 15:00-04.14 17:00"
$ws.Range("F7").Value = 240
$ws.Range("E8").Value = "2024.04.14 19:00-04.14 20:30"
$ws.Range("E9").Value = "2024.04.19 19:30-04.19 21:00"
$ws.Range("E10").Value = "2024.04.24 20:00-04.24 21:30"
$ws.Range("E11").Value = "2024.04.28 19:30-04.28 21:30"
$ws.Range("E12").Value = "2024.04.28 19:00-04.28 20:30"

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("E2").Value = "2024.02.01 00:00-03.17 23:59"
$ws.Range("E3").Value = "2024.03.21 00:00-04.28 23:59"
$ws.Range("F3").Value = 101

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("E2").Value = "2024.02.01 00:00-03.17 23:59"
$ws.Range("E3").Value = "2024.03.16 13:00-03.17 19:00"
$ws.Range("F3").Value = 889
$ws.Range("E4").Value = "2024.03.16 10:00-03.16 21:00"
$ws.Range("E5").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("F5").Value = 1136
$ws.Range("E6").Value = "2024.03.17 19:00-03.17 20:30"
$ws.Range("E7").Value = "2024.03.17 19:30-03.17 21:00"
$ws.Range("F7").Value = 34
$ws.Range("E8").Value = "2024.03.21 00:00-04.28 23:59"
$ws.Range("F8").Value = 101
$ws.Range("E9").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F9").Value = 535
$ws.Range("E10").Value = "2024.03.23 10:00-03.24 17:00"
$ws.Range("F10").Value = 230
$ws.Range("E11").Value = "2024.03.23 13:00-03.23 19:00"
$ws.Range("E12").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F12").Value = 692
$ws.Range("E13").Value = "2024.03.23 20:00-03.23 22:00"
$ws.Range("E14").Value = "2024.03.23 09:30-03.24 18:00"
$ws.Range("F14").Value = 265
$ws.Range("E15").Value = "2024.03.23 09:30-03.23 18:00"
$ws.Range("E16").Value = "2024.03.24 09:30-03.24 18:00"
$ws.Range("F16").Value = 101
$ws.Range("E17").Value = "2024.03.24 10:00-03.24 17:30"
$ws.Range("F17").Value = 228
$ws.Range("E18").Value = "2024.03.24 11:00-03.24 17:00"
$ws.Range("F18").Value = 165
$ws.Range("E19").Value = "2024.03.30 09:00-03.31 17:00"
$ws.Range("F19").Value = 3349
$ws.Range("E20").Value = "2024.03.30 14:00-03.30 18:00"
$ws.Range("F20").Value = 14
$ws.Range("E21").Value = "2024.03.30 10:00-03.30 17:00"
$ws.Range("F21").Value = 14
$ws.Range("E22").Value = "2024.03.31 19:00-03.31 20:30"
$ws.Range("F22").Value = 243
$ws.Range("E23").Value = "2024.04.04 10:00-04.05 17:30"
$ws.Range("E24").Value = "2024.04.05 11:30-04.06 18:00"
$ws.Range("F24").Value = 42
$ws.Range("E25").Value = "2024.04.06 10:00-04.06 17:00"
$ws.Range("E26").Value = "2024.04.06 10:00-04.06 17:00"
$ws.Range("F26").Value = 290
$ws.Range("E27").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E28").Value = "2024.04.13 19:00-04.13 22:00"
$ws.Range("E29").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E30").Value = "2024.04.14 10:00-04.14 17:00"
$ws.Range("F30").Value = 8
$ws.Range("E31").Value = "2024.04.14 15:00-04.14 17:00"
$ws.Range("F31").Value = 240
$ws.Range("E32").Value = "2024.04.14 19:00-04.14 20:30"
$ws.Range("E33").Value = "2024.04.19 19:30-04.19 21:00"
$ws.Range("E34").Value = "2024.04.20 10:00-04.20 17:00"
$ws.Range("F34").Value = 681
$ws.Range("E35").Value = "2024.04.24 20:00-04.24 21:30"
$ws.Range("E36").Value = "2024.04.28 19:30-04.28 21:30"
$ws.Range("E37").Value = "2024.04.28 19:00-04.28 20:30"
$ws.Range("E38").Value = "2024.05.04 10:00-05.04 17:00"
$ws.Range("F38").Value = 62
$ws.Range("E39").Value = "2024.05.04 09:30-05.04 16:00"
$ws.Range("F39").Value = 258
$ws.Range("E40").Value = "2024.05.05 10:00-05.05 17:00"
$ws.Range("F40").Value = 972
$ws.Range("E41").Value = "2024.05.10 10:30-05.12 18:30"
$ws.Range("F41").Value = 56
$ws.Range("E42").Value = "2024.05.18 10:00-05.18 17:00"
$ws.Range("F42").Value = 1614
$ws.Range("E43").Value = "2024.05.18 10:00-05.18 17:00"
$ws.Range("F43").Value = 351
$ws.Range("E44").Value = "2024.05.25 10:30-05.25 17:00"
